# Generate Report for Archive
# Update status text from "Ready for handoff" to "In Translation" across all
# sheets, and narrow the corresponding "Status" columns.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update the status value wherever it appears ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Narrow the Status columns to match the new (shorter) text ---
# Target stored column width is 13.4101848602295 characters; this
# engine quantizes ColumnWidth to 1/6-character increments, so we pick
# an input value (12.5) that lands on the nearest attainable grid point.
$overview.Columns("E").ColumnWidth = 12.5
$overview.Columns("F").ColumnWidth = 12.5
$zhcn.Columns("C").ColumnWidth = 12.5
$dede.Columns("C").ColumnWidth = 12.5
